$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il16"
$ws.Range("C2").Value = "Kcnj15"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.7754575
$ws.Range("H2").Value = 7.550915
$ws.Range("I2").Value = 0.1708449704380497
$ws.Range("J2").Value = 0.1225509496985904
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.1986275
$ws.Range("N2").Value = 4.397254999999999
$ws.Range("O2").Value = 0.4935104206364654
$ws.Range("P2").Value = 0.3937861339617307
$ws.Range("Q2").Value = 8.300824684581249
$ws.Range("R2").Value = 33.203298738325
$ws.Range("S2").Value = 0.08431377322450638
$ws.Range("T2").Value = 0.04825886469514645

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il16"
$ws.Range("C3").Value = "Kcnj15"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.7754575
$ws.Range("H3").Value = 7.550915
$ws.Range("I3").Value = 0.1708449704380497
$ws.Range("J3").Value = 0.1225509496985904
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.173422
$ws.Range("N3").Value = 0.520266
$ws.Range("O3").Value = 0.0389268141909519
$ws.Range("P3").Value = 0.0465912340248027
$ws.Range("Q3").Value = 0.654747390565
$ws.Range("R3").Value = 3.92848434339
$ws.Range("S3").Value = 0.006650450419700631
$ws.Range("T3").Value = 0.00570979997736885

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il16"
$ws.Range("C4").Value = "Kcnj15"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.7754575
$ws.Range("H4").Value = 7.550915
$ws.Range("I4").Value = 0.1708449704380497
$ws.Range("J4").Value = 0.1225509496985904
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.083028666666667
$ws.Range("N4").Value = 6.249086
$ws.Range("O4").Value = 0.4675627651725827
$ws.Range("P4").Value = 0.5596226320134666
$ws.Range("Q4").Value = 7.864386202281668
$ws.Range("R4").Value = 47.18631721369
$ws.Range("S4").Value = 0.07988074679384265
$ws.Range("T4").Value = 0.06858228502607512

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il16"
$ws.Range("C5").Value = "Kcnj15"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.327317666666667
$ws.Range("H5").Value = 3.981953
$ws.Range("I5").Value = 0.06006306454875145
$ws.Range("J5").Value = 0.06462688585491311
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.1986275
$ws.Range("N5").Value = 4.397254999999999
$ws.Range("O5").Value = 0.4935104206364654
$ws.Range("P5").Value = 0.3937861339617307
$ws.Range("Q5").Value = 2.918277123169167
$ws.Range("R5").Value = 17.509662739015
$ws.Range("S5").Value = 0.0296417482501695
$ws.Range("T5").Value = 0.0254491715307923

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il16"
$ws.Range("C6").Value = "Kcnj15"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.327317666666667
$ws.Range("H6").Value = 3.981953
$ws.Range("I6").Value = 0.06006306454875145
$ws.Range("J6").Value = 0.06462688585491311
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.173422
$ws.Range("N6").Value = 0.520266
$ws.Range("O6").Value = 0.0389268141909519
$ws.Range("P6").Value = 0.0465912340248027
$ws.Range("Q6").Value = 0.2301860843886667
$ws.Range("R6").Value = 2.071674759498
$ws.Range("S6").Value = 0.002338063753428398
$ws.Range("T6").Value = 0.003011046363160468

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il16"
$ws.Range("C7").Value = "Kcnj15"
$ws.Range("D7").Value = "Neutro"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.327317666666667
$ws.Range("H7").Value = 3.981953
$ws.Range("I7").Value = 0.06006306454875145
$ws.Range("J7").Value = 0.06462688585491311
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.083028666666667
$ws.Range("N7").Value = 6.249086
$ws.Range("O7").Value = 0.4675627651725827
$ws.Range("P7").Value = 0.5596226320134666
$ws.Range("Q7").Value = 2.764840749439778
$ws.Range("R7").Value = 24.883566744958
$ws.Range("S7").Value = 0.02808325254515355
$ws.Range("T7").Value = 0.03616666796096035

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Il16"
$ws.Range("C8").Value = "Kcnj15"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.188097999999999
$ws.Range("H8").Value = 15.564294
$ws.Range("I8").Value = 0.2347690179109961
$ws.Range("J8").Value = 0.2526076655727249
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.1986275
$ws.Range("N8").Value = 4.397254999999999
$ws.Range("O8").Value = 0.4935104206364654
$ws.Range("P8").Value = 0.3937861339617307
$ws.Range("Q8").Value = 11.406694935495
$ws.Range("R8").Value = 68.44016961296998
$ws.Range("S8").Value = 0.1158609567816656
$ws.Range("T8").Value = 0.09947339603498112

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Il16"
$ws.Range("C9").Value = "Kcnj15"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.188097999999999
$ws.Range("H9").Value = 15.564294
$ws.Range("I9").Value = 0.2347690179109961
$ws.Range("J9").Value = 0.2526076655727249
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.173422
$ws.Range("N9").Value = 0.520266
$ws.Range("O9").Value = 0.0389268141909519
$ws.Range("P9").Value = 0.0465912340248027
$ws.Range("Q9").Value = 0.8997303313559998
$ws.Range("R9").Value = 8.097572982204
$ws.Range("S9").Value = 0.009138809938013604
$ws.Range("T9").Value = 0.01176930286315792

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Il16"
$ws.Range("C10").Value = "Kcnj15"
$ws.Range("D10").Value = "Neutro"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.188097999999999
$ws.Range("H10").Value = 15.564294
$ws.Range("I10").Value = 0.2347690179109961
$ws.Range("J10").Value = 0.2526076655727249
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.083028666666667
$ws.Range("N10").Value = 6.249086
$ws.Range("O10").Value = 0.4675627651725827
$ws.Range("P10").Value = 0.5596226320134666
$ws.Range("Q10").Value = 10.806956859476
$ws.Range("R10").Value = 97.262611735284
$ws.Range("S10").Value = 0.1097692511913169
$ws.Range("T10").Value = 0.1413649666745859

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Il16"
$ws.Range("C11").Value = "Kcnj15"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.677527666666667
$ws.Range("H11").Value = 23.032583
$ws.Range("I11").Value = 0.3474193491117236
$ws.Range("J11").Value = 0.3738175996765437
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.1986275
$ws.Range("N11").Value = 4.397254999999999
$ws.Range("O11").Value = 0.4935104206364654
$ws.Range("P11").Value = 0.3937861339617307
$ws.Range("Q11").Value = 16.88002345994417
$ws.Range("R11").Value = 101.280140759665
$ws.Range("S11").Value = 0.1714550691173737
$ws.Range("T11").Value = 0.14720418738348

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Il16"
$ws.Range("C12").Value = "Kcnj15"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.677527666666667
$ws.Range("H12").Value = 23.032583
$ws.Range("I12").Value = 0.3474193491117236
$ws.Range("J12").Value = 0.3738175996765437
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.173422
$ws.Range("N12").Value = 0.520266
$ws.Range("O12").Value = 0.0389268141909519
$ws.Range("P12").Value = 0.0465912340248027
$ws.Range("Q12").Value = 1.331452203008667
$ws.Range("R12").Value = 11.983069827078
$ws.Range("S12").Value = 0.01352392844921352
$ws.Range("T12").Value = 0.01741662326911986

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Il16"
$ws.Range("C13").Value = "Kcnj15"
$ws.Range("D13").Value = "Neutro"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.677527666666667
$ws.Range("H13").Value = 23.032583
$ws.Range("I13").Value = 0.3474193491117236
$ws.Range("J13").Value = 0.3738175996765437
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.083028666666667
$ws.Range("N13").Value = 6.249086
$ws.Range("O13").Value = 0.4675627651725827
$ws.Range("P13").Value = 0.5596226320134666
$ws.Range("Q13").Value = 15.99251021879311
$ws.Range("R13").Value = 143.932591969138
$ws.Range("S13").Value = 0.1624403515451364
$ws.Range("T13").Value = 0.2091967890239438

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Il16"
$ws.Range("C14").Value = "Kcnj15"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.224085333333333
$ws.Range("H14").Value = 9.672256
$ws.Range("I14").Value = 0.145894573959072
$ws.Range("J14").Value = 0.1569802015421825
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.1986275
$ws.Range("N14").Value = 4.397254999999999
$ws.Range("O14").Value = 0.4935104206364654
$ws.Range("P14").Value = 0.3937861339617307
$ws.Range("Q14").Value = 7.088562676213333
$ws.Range("R14").Value = 42.53137605728
$ws.Range("S14").Value = 0.07200049256311951
$ws.Range("T14").Value = 0.06181662667382938

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Il16"
$ws.Range("C15").Value = "Kcnj15"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.224085333333333
$ws.Range("H15").Value = 9.672256
$ws.Range("I15").Value = 0.145894573959072
$ws.Range("J15").Value = 0.1569802015421825
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.173422
$ws.Range("N15").Value = 0.520266
$ws.Range("O15").Value = 0.0389268141909519
$ws.Range("P15").Value = 0.0465912340248027
$ws.Range("Q15").Value = 0.5591273266773333
$ws.Range("R15").Value = 5.032145940096001
$ws.Range("S15").Value = 0.005679210971972885
$ws.Range("T15").Value = 0.007313901307312521

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Il16"
$ws.Range("C16").Value = "Kcnj15"
$ws.Range("D16").Value = "Neutro"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.224085333333333
$ws.Range("H16").Value = 9.672256
$ws.Range("I16").Value = 0.145894573959072
$ws.Range("J16").Value = 0.1569802015421825
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.083028666666667
$ws.Range("N16").Value = 6.249086
$ws.Range("O16").Value = 0.4675627651725827
$ws.Range("P16").Value = 0.5596226320134666
$ws.Range("Q16").Value = 6.71586217311289
$ws.Range("R16").Value = 60.44275955801601
$ws.Range("S16").Value = 0.06821487042397957
$ws.Range("T16").Value = 0.08784967356104065

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Il16"
$ws.Range("C17").Value = "Kcnj15"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.9062475
$ws.Range("H17").Value = 1.812495
$ws.Range("I17").Value = 0.04100902403140717
$ws.Range("J17").Value = 0.02941669765504533
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.1986275
$ws.Range("N17").Value = 4.397254999999999
$ws.Range("O17").Value = 0.4935104206364654
$ws.Range("P17").Value = 0.3937861339617307
$ws.Range("Q17").Value = 1.99250067530625
$ws.Range("R17").Value = 7.970002701224999
$ws.Range("S17").Value = 0.02023838069963067
$ws.Range("T17").Value = 0.01158388764350141

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Il16"
$ws.Range("C18").Value = "Kcnj15"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.9062475
$ws.Range("H18").Value = 1.812495
$ws.Range("I18").Value = 0.04100902403140717
$ws.Range("J18").Value = 0.02941669765504533
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.173422
$ws.Range("N18").Value = 0.520266
$ws.Range("O18").Value = 0.0389268141909519
$ws.Range("P18").Value = 0.0465912340248027
$ws.Range("Q18").Value = 0.157163253945
$ws.Range("R18").Value = 0.94297952367
$ws.Range("S18").Value = 0.001596350658622868
$ws.Range("T18").Value = 0.001370560244683082

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Il16"
$ws.Range("C19").Value = "Kcnj15"
$ws.Range("D19").Value = "Neutro"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.9062475
$ws.Range("H19").Value = 1.812495
$ws.Range("I19").Value = 0.04100902403140717
$ws.Range("J19").Value = 0.02941669765504533
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 2.083028666666667
$ws.Range("N19").Value = 6.249086
$ws.Range("O19").Value = 0.4675627651725827
$ws.Range("P19").Value = 0.5596226320134666
$ws.Range("Q19").Value = 1.887739521595
$ws.Range("R19").Value = 11.32643712957
$ws.Range("S19").Value = 0.01917429267315363
$ws.Range("T19").Value = 0.01646224976686084
